$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the existing "derivedVoltageCreationServiceUrl" row down from row 7 to
# row 8, insert a new "derivedVdiCreationServiceUrl" row at 9, and retarget
# row 7 to a new "derivedFrequencyCreationServiceUrl" entry. Set cell values
# in this exact order so the shared-string table regenerates with the same
# index assignments as the target workbook (idx 8 reused, idx 9 = Vdi,
# idx 10 = Frequency).
$ws.Range("A8").Value = "derivedVoltageCreationServiceUrl"
$ws.Range("A9").Value = "derivedVdiCreationServiceUrl"
$ws.Range("A7").Value = "derivedFrequencyCreationServiceUrl"

# B8 / B9 get the same "http://google.com" hyperlink + Hyperlink style as the
# rest of column B.
$ws.Range("B8").Value = "http://google.com"
$ws.Range("B9").Value = "http://google.com"

$ws.Hyperlinks.Add($ws.Range("B8"), "http://google.com/")
$ws.Hyperlinks.Add($ws.Range("B9"), "http://google.com/")

# Hyperlinks.Add() stamps its own style on the cell; reapply the existing
# Hyperlink cell style (copied from B3) so B8/B9 match B3:B7 exactly.
$ws.Range("B8").Style = $ws.Range("B3").Style
$ws.Range("B9").Style = $ws.Range("B3").Style
